# Updates crypto price/volume data per upstream GitHub Actions refresh (commit: "Updated cryptos list on Sun Sep  8 19:27:06 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.288.45"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.263.97"
$ws.Range("E3").Value = "  -1.01%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'495.64"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'128.71"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.525"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").Value = "'0.0954"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.42%  "
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("D11").Value = "'0.337"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.82%  "
$ws.Range("E12").Value = "  +3.05%  "
$ws.Range("D13").Value = "'22.98"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.24%  "
$ws.Range("D14").Value = "2.666.74"
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "54.254.45"
$ws.Range("E15").Value = "  -0.04%  "
$ws.Range("E16").Value = "  +0.20%  "
$ws.Range("D17").Value = "2.279.61"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'10.23"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.83%  "
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("D20").Value = "'302.97"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("D21").Value = "'6.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.03%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'60.73"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.09%  "
$ws.Range("E24").Value = "  -1.54%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "'7.31"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("D27").Value = "'171.39"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("B28").Value = "PancakeSwap"
$ws.Range("C28").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D28").Value = "'1.61"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").Value = "'5.97"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").Value = "0.0₃0688"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +1.06%  "
$ws.Range("D33").Value = "'17.78"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "'0.937"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.54%  "
$ws.Range("D36").Value = "'1.19"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.14%  "
$ws.Range("D37").Value = "'3.70"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "'0.374"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "'1.39"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "'3.37"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.15%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").Value = "'124.82"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'4.80"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.92%  "
$ws.Range("E43").Value = "  +2.15%  "
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("E45").Value = "  +0.64%  "
$ws.Range("D46").Value = "'241.20"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.08%  "
$ws.Range("D47").Value = "'0.373"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.09%  "
$ws.Range("D48").Value = "'0.0204"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").Value = "'10.82"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.90%  "
$ws.Range("D50").Value = "'16.08"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("E51").Value = "  -0.46%  "
